# lugares_visitar.xlsx - mark a few places as visited and add newly visited places
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Mark previously-unvisited establishments as visited ("N" -> "S")
$ws.Range("C11").Value = "S"   # Madero
$ws.Range("C13").Value = "S"   # YoouGeek
$ws.Range("C20").Value = "S"   # Home BBQ

# Append newly discovered establishments to the bottom of the list
$ws.Range("A26").Value = "Don Hamburgo"
$ws.Range("B26").Value = "Campinas"
$ws.Range("C26").Value = "N"

$ws.Range("A27").Value = "L'Entrecôte de Paris"
$ws.Range("B27").Value = "Campinas"
$ws.Range("C27").Value = "N"

$ws.Range("A28").Value = "Bistrô Paris 6"
$ws.Range("B28").Value = "Campinas"
$ws.Range("C28").Value = "N"

$ws.Range("A29").Value = "Beco Hexagonal"
$ws.Range("B29").Value = "São Paulo"
$ws.Range("C29").Value = "N"

# Keep selection/view consistent with the saved file
$ws.Range("C21").Select() | Out-Null
